$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data: refresh the dSF (F) column values for a set of rows
# (mean calculation adjustments) -- F column recomputed from updated source data.
$ws.Range("F8").Value  = 1
$ws.Range("F9").Value  = 0
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = -4
$ws.Range("F27").Value = 3
$ws.Range("F29").Value = -2
$ws.Range("F31").Value = 1
$ws.Range("F35").Value = -1
$ws.Range("F46").Value = 1
$ws.Range("F53").Value = -3
$ws.Range("F57").Value = 2
$ws.Range("F60").Value = 1
